# OW-535: update the bilateral trade portfolio to match the acuo-data test branch.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Position Account ID (column B, row 2): swap in the acuo-data test account id.
$ws.Range("B2").Value = "ACUOSG8745"

# Portfolio ID (column AP, row 2): swap in the acuo-data test portfolio id.
$ws.Range("AP2").Value = "p1"

# Re-anchor the sheet selection on the single cell A2 (was a full-row range A2:AR2).
$ws.Range("A2").Select()
